# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# New F-column values, keyed by row number. Both the "展览" sheet and the
# "全部类型" sheet carry identical copies of this data and both need the
# same update.
$updates = @{
    3  = 1099
    8  = 11306
    9  = 4297
    13 = 2515
    14 = 1074
    15 = 113
    18 = 494
    19 = 11258
    20 = 11116
    22 = 38
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
